$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample")

# --- Row 2: B2 and D2 swap content (CODE <-> DESCRIPTIVE_TEXT) ---
$ws.Range("B2").Value = "CODE"
$ws.Range("D2").Value = "DESCRIPTIVE_TEXT"

# --- Row 4: drop the ":CODE" / ":DESCRIPTIVE_TEXT" suffixes ---
$ws.Range("A4").Value = "A2"
$ws.Range("B4").Value = "B2"

# --- Row 5 takes on what used to be row 6's content; alignment goes right->left ---
$ws.Range("A5").Value = "A4"
$ws.Range("B5").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("B5").Value = "X; Y; Z"

# --- Row 6 takes on what used to be row 7's content; B6 loses the border/font accent ---
$ws.Range("A6").Value = "A5; A6; A7"
$ws.Range("C6").Copy()
$ws.Range("B6").PasteSpecial(-4122)            # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B6").Value = "B5; B6; B7"

# --- Row 7 becomes blank, but keeps (gains) a quote-prefix style ---
$ws.Range("A7").Value = "'"
$ws.Range("B7").Value = "'"
